# Standard Build Flow.docx — add clause for redirecting old pages
#
# This script applies five related changes:
#   1. Insert a new bulleted paragraph (numId=2) right before the
#      "You are finished cropping..." paragraph, moving the
#      <w:lastRenderedPageBreak/> marker onto the new paragraph's run.
#   2. Remove the (hidden) _GoBack bookmark from the "Create the city
#      pages." paragraph.
#   3. Add a <w:lastRenderedPageBreak/> to the "Find the newest server..."
#      run.
#   4. Remove the <w:lastRenderedPageBreak/> that used to sit on the
#      "FTP into the newly created..." run.
#   5. Insert a new bulleted paragraph (numId=4) right after the
#      "Update the sitemap..." paragraph, carrying the _GoBack bookmark
#      that used to live on "Create the city pages.".

$d = $word.ActiveDocument

function Find-ParagraphByText($doc, [string]$pattern) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text -match $pattern) {
            return $p
        }
    }
    return $null
}

function Remove-EmptyParagraphMarker($doc) {
    # InsertXML of "<w:p>...</w:p><w:p/>" leaves a stray empty paragraph
    # behind as a side effect of forcing the paragraph split; clean it up.
    foreach ($p in $doc.Paragraphs) {
        $t = $p.Range.Text
        if ($t -eq "" -or $t -eq "`r") {
            $p.Range.Delete()
            return
        }
    }
}

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- 1. New clause about redirecting an existing live website -------------
$finishedPara = Find-ParagraphByText $d "You are finished cropping"
$insertPoint = $d.Range($finishedPara.Range.Start, $finishedPara.Range.Start)
$newParaXml = $pkgHeader + `
    '<w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>If they have an existing website that is live, perform rewrites or redirects with a 301 HTTP code to any corresponding new URLs. The new URLs are typically written at the top of the writer-provided content for pages.</w:t></w:r></w:p><w:p/></w:body>' + `
    $pkgFooter
$insertPoint.InsertXML($newParaXml)
Remove-EmptyParagraphMarker $d

# The page break marker used to live on the "You are finished cropping..."
# run; now that it has moved to the newly-inserted paragraph, rewrite that
# run without it.
$finishedPara = Find-ParagraphByText $d "You are finished cropping"
$finishedXml = $pkgHeader + `
    '<w:body><w:p><w:r><w:t>You are finished cropping the template pieces of the layout and setting up the website for the rest of the build.</w:t></w:r></w:p></w:body>' + `
    $pkgFooter
$finishedPara.Range.InsertXML($finishedXml)

# --- 2. Drop the _GoBack bookmark from "Create the city pages." -----------
$cityPara = Find-ParagraphByText $d "Create the city pages"
$cityXml = $pkgHeader + `
    '<w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>Create the city pages.</w:t></w:r></w:p></w:body>' + `
    $pkgFooter
$cityPara.Range.InsertXML($cityXml)

# --- 3/4. Move <w:lastRenderedPageBreak/> from "FTP into the newly
#          created..." onto "Find the newest server...".
$findServerPara = Find-ParagraphByText $d "Find the newest server in the background"
$findServerXml = $pkgHeader + `
    '<w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Find the newest server in the background and access the WHM. Login.</w:t></w:r></w:p></w:body>' + `
    $pkgFooter
$findServerPara.Range.InsertXML($findServerXml)

$ftpPara = Find-ParagraphByText $d "FTP into the newly created"
$ftpXml = $pkgHeader + `
    '<w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>FTP into the newly created. If there are any files in here that aren' + [char]0x2019 + 't defaults (possibly because the client provided an already made server). Back these files up.</w:t></w:r><w:r><w:t xml:space="preserve"> Upload the new website files.</w:t></w:r></w:p></w:body>' + `
    $pkgFooter
$ftpPara.Range.InsertXML($ftpXml)

# --- 5. New clause verifying redirects are in place; carries the _GoBack
#        bookmark that used to sit on "Create the city pages.".
$sitemapPara = Find-ParagraphByText $d "Update the sitemap full URL"
$afterSitemap = $d.Range($sitemapPara.Range.End, $sitemapPara.Range.End)
$verifyXml = $pkgHeader + `
    '<w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>Verify that, if they had a previous website, any redirects from the old URLs to corresponding new URLs are in place.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p/></w:body>' + `
    $pkgFooter
$afterSitemap.InsertXML($verifyXml)
Remove-EmptyParagraphMarker $d

Write-Output "edit complete"
